$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 now documents the "add line drawing demo to appendix" follow-up under Scala.React
$ws.Range("A9").Value = "Scala.React"
$ws.Range("B9").Value = "add line drawing demo to appendix"

# Two new Scala.Swing rows
$ws.Range("A10").Value = "Scala.Swing"
$ws.Range("B10").Value = "example"

$ws.Range("A11").Value = "Scala.Swing"
$ws.Range("B11").Value = "explain event mechanism"

# The old "Scomm / explain a copy/new folder" row moves down and gains a status
$ws.Range("A12").Value = "Scomm"
$ws.Range("B12").Value = "explain a copy/new folder"
$ws.Range("C12").Value = "IN PROGRESS"

# The old "Scomm / case study" row moves down too
$ws.Range("A13").Value = "Scomm"
$ws.Range("B13").Value = "case study"

# Grow the table to cover the newly added rows
$lo = $ws.ListObjects.Item(1)
$null = $lo.Resize($ws.Range("A1:C13"))

# Column C now holds the status text; size it to fit like the source workbook
$ws.Columns.Item(3).ColumnWidth = 11

# Leave the selection where the author finished editing
$null = $ws.Range("B13").Select()
